$d = $word.ActiveDocument
$p = $d.Paragraphs(16)
$r = $p.Range
$r.Collapse(1)  # wdCollapseStart = 1
$r.InsertParagraphBefore()

$newPara = $d.Paragraphs(16)
$pStart = $newPara.Range.Start

# Insert the full text first (bold, from template)
$ip = $d.Range($pStart, $pStart)
$fullText = "ALIGN ALT/GPS"
$ip.InsertBefore($fullText)

$afterBoldPos = $pStart + $fullText.Length
$ip2 = $d.Range($afterBoldPos, $afterBoldPos)
$ip2.InsertBefore(":")
$afterColonPos = $afterBoldPos + 1
$ip3 = $d.Range($afterColonPos, $afterColonPos)
$ip3.InsertBefore(" It should be possible to overlay interval data by aligning to GPS data or the altitude profile if not then use time/distance only.")

Write-Output $d.Paragraphs(16).Range.Text
Write-Output "---"
# Now set bold=false on range from afterBoldPos to new end (excluding para mark)
$paraEnd = $d.Paragraphs(16).Range.End

$colonRange = $d.Range($afterBoldPos, $afterColonPos)
$colonRange.Font.Bold = 0
$colonRange.ClearFormatting()
$tailRange = $d.Range($afterColonPos, $paraEnd - 1)
$tailRange.Font.Bold = 0
$tailRange.ClearFormatting()
Write-Output $d.Paragraphs(16).Range.Text
